$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.483.09'
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').Value = '1.901.01'
$ws.Range('E3').Value = '  +1.45%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.24%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.05'
$ws.Range('E5').Value = '  +1.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  +0.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4908'
$ws.Range('E7').Value = '  +0.88%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2934'
$ws.Range('E8').Value = '  +1.29%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06683'
$ws.Range('E9').Value = '  +0.31%  '

$ws.Range('D10').Value = '1.909.72'
$ws.Range('E10').Value = '  +2.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.02'
$ws.Range('E11').Value = '  +2.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07336'
$ws.Range('E12').Value = '  +1.59%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.179'
$ws.Range('E13').Value = '  +3.65%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.26'
$ws.Range('E14').Value = '  -0.89%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6683'
$ws.Range('E15').Value = '  +2.47%  '

$ws.Range('D16').Value = '30.462.95'
$ws.Range('E16').Value = '  +0.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007887'
$ws.Range('E17').Value = '  +0.44%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.45'
$ws.Range('E18').Value = '  +3.45%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9997'
$ws.Range('E19').Value = '  +0.12%  '

$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.146.77'
$ws.Range('E20').Value = '  +1.54%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.422'
$ws.Range('E21').Value = '  +15.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9989'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '199.74'
$ws.Range('E23').Value = '  -7.43%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.138'
$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.506'
$ws.Range('E25').Value = '  +1.56%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.29'
$ws.Range('E26').Value = '  +3.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.47'
$ws.Range('E27').Value = '  -1.85%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.943'
$ws.Range('E28').Value = '  +6.59%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.483'
$ws.Range('E29').Value = '  +5.58%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.344'
$ws.Range('E30').Value = '  +2.17%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09166'
$ws.Range('E31').Value = '  +1.64%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.132'
$ws.Range('E32').Value = '  +5.47%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05158'
$ws.Range('E33').Value = '  +1.18%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7414'
$ws.Range('E34').Value = '  +2.85%  '

$ws.Range('E35').Value = '  +2.98%  '

$ws.Range('E36').Value = '  +1.62%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01840'
$ws.Range('E37').Value = '  +1.82%  '

$ws.Range('E38').Value = '  +0.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9266'
$ws.Range('E39').Value = '  +0.89%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.069'
$ws.Range('E40').Value = '  +1.49%  '

$ws.Range('E41').Value = '  +0.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.92'
$ws.Range('E42').Value = '  +2.46%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.927'
$ws.Range('E43').Value = '  +3.47%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.65'
$ws.Range('E44').Value = '  +22.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9954'
$ws.Range('E45').Value = '  +0.22%  '

$ws.Range('E46').Value = '  +3.64%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.599'
$ws.Range('E47').Value = '  +3.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.014'
$ws.Range('E48').Value = '  +4.72%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.11'
$ws.Range('E49').Value = '  +6.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05839'
$ws.Range('E50').Value = '  +0.43%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3945'
$ws.Range('E51').Value = '  -1.48%  '
